# Update countries & provincias Spain
#
# 1) Re-order three pairs of countries (the row that used to be first now
#    comes second, and vice-versa), while each country keeps the COVID
#    figures that belong to it.
# 2) Refresh the numeric COVID figures for a number of countries.
# 3) Bump the "Datos actualizados..." timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the country names so that the data moves down one row -------
# Serbia (row 61) <-> Moldavia (row 62)
$ws.Range("A61").Value = "Moldavia"
$ws.Range("A62").Value = "Serbia"

# Etiopia (row 80) <-> Republica de Macedonia (row 81)
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("A81").Value = "Etiopia"

# Grecia (row 101) <-> Albania (row 102)
$ws.Range("A101").Value = "Albania"
$ws.Range("A102").Value = "Grecia"

# --- 2) Refresh the numeric figures ---------------------------------------
# Estados Unidos (row 4)
$ws.Range("B4").Value = 3795389
$ws.Range("C4").Value = 25377
$ws.Range("D4").Value = 1756493
$ws.Range("E4").Value = 1896465
$ws.Range("G4").Value = 367
$ws.Range("H4").Value = 142431

# Brasil (row 5)
$ws.Range("B5").Value = 2053174
$ws.Range("C5").Value = 4477
$ws.Range("E5").Value = 608302
$ws.Range("G5").Value = 165
$ws.Range("H5").Value = 78097

# Chile (row 11)
$ws.Range("B11").Value = 328846
$ws.Range("C11").Value = 2307
$ws.Range("D11").Value = 299449
$ws.Range("E11").Value = 20952
$ws.Range("G11").Value = 98
$ws.Range("H11").Value = 8445

# Canada (row 24)
$ws.Range("B24").Value = 109835
$ws.Range("C24").Value = 166
$ws.Range("D24").Value = 96819
$ws.Range("E24").Value = 4175
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 8841

# Argelia (row 60)
$ws.Range("B60").Value = 22549
$ws.Range("C60").Value = 601
$ws.Range("E60").Value = 6051
$ws.Range("G60").Value = 11
$ws.Range("H60").Value = 1068

# Serbia, now at row 61 (new figures)
$ws.Range("B61").Value = 20794
$ws.Range("C61").Value = 300
$ws.Range("D61").Value = 14183
$ws.Range("E61").Value = 5931
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 680

# Moldavia, now at row 62 (figures carried over from the old row 61)
$ws.Range("B62").Value = 20498
$ws.Range("C62").Value = 389
$ws.Range("D62").Value = 14047
$ws.Range("E62").Value = 5990
$ws.Range("G62").Value = 9
$ws.Range("H62").Value = 461

# Chequia (row 68)
$ws.Range("B68").Value = 13795
$ws.Range("C68").Value = 53
$ws.Range("E68").Value = 4712

# Etiopia, now at row 80 (new figures)
$ws.Range("B80").Value = 9026
$ws.Range("C80").Value = 240
$ws.Range("D80").Value = 4727
$ws.Range("E80").Value = 3885
$ws.Range("G80").Value = 8
$ws.Range("H80").Value = 414

# Republica de Macedonia, now at row 81 (figures carried over from the old row 80)
$ws.Range("B81").Value = 8803
$ws.Range("D81").Value = 2430
$ws.Range("E81").Value = 6223
$ws.Range("H81").Value = 150

# Luxemburgo (row 96)
$ws.Range("B96").Value = 5483
$ws.Range("C96").Value = 74
$ws.Range("E96").Value = 1039

# Grecia, now at row 101 (new figures)
$ws.Range("B101").Value = 4008
$ws.Range("C101").Value = 102
$ws.Range("D101").Value = 2264
$ws.Range("E101").Value = 1633
$ws.Range("G101").Value = 4
$ws.Range("H101").Value = 111

# Albania, now at row 102 (figures carried over from the old row 101)
$ws.Range("B102").Value = 3983
$ws.Range("C102").Value = 19
$ws.Range("D102").Value = 1374
$ws.Range("E102").Value = 2415
$ws.Range("H102").Value = 194

# --- 3) Bump the "updated at" timestamp -----------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 19:09"
